$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 401 (existing rows 401:510 shift down to 404:513).
$ws.Rows("401:403").Insert()

# Row 401 (new) - "Especial" quality, new week's price data.
$ws.Range("A401").Value = 8
$ws.Range("B401").Value = "Terminal La Palmera de La Serena"
$ws.Range("C401").Value = "Coquimbo"
$ws.Range("D401").Value = 44543
$ws.Range("E401").Value = 4
$ws.Range("F401").Value = "Fruta"
$ws.Range("G401").Value = 100101
$ws.Range("H401").Value = "Berries"
$ws.Range("I401").Value = 100112025
$ws.Range("J401").Value = "Frutilla"
$ws.Range("K401").Value = "Sin especificar"
$ws.Range("L401").Value = "Especial"
$ws.Range("M401").Value = 300
$ws.Range("N401").Value = 11500
$ws.Range("O401").Value = 12000
$ws.Range("P401").Value = 11750
$ws.Range("Q401").Value = "$/bandeja 7 kilos"
$ws.Range("R401").Value = "Provincia de Melipilla"
$ws.Range("S401").Value = 1679
$ws.Range("T401").Value = 7

# Row 402 (new) - "Primera" quality, new week's price data.
$ws.Range("A402").Value = 8
$ws.Range("B402").Value = "Terminal La Palmera de La Serena"
$ws.Range("C402").Value = "Coquimbo"
$ws.Range("D402").Value = 44543
$ws.Range("E402").Value = 4
$ws.Range("F402").Value = "Fruta"
$ws.Range("G402").Value = 100101
$ws.Range("H402").Value = "Berries"
$ws.Range("I402").Value = 100112025
$ws.Range("J402").Value = "Frutilla"
$ws.Range("K402").Value = "Sin especificar"
$ws.Range("L402").Value = "Primera"
$ws.Range("M402").Value = 300
$ws.Range("N402").Value = 9500
$ws.Range("O402").Value = 10000
$ws.Range("P402").Value = 9750
$ws.Range("Q402").Value = "$/bandeja 7 kilos"
$ws.Range("R402").Value = "Provincia de Melipilla"
$ws.Range("S402").Value = 1393
$ws.Range("T402").Value = 7

# Row 403 (new) - "Segunda" quality, new week's price data.
$ws.Range("A403").Value = 8
$ws.Range("B403").Value = "Terminal La Palmera de La Serena"
$ws.Range("C403").Value = "Coquimbo"
$ws.Range("D403").Value = 44543
$ws.Range("E403").Value = 4
$ws.Range("F403").Value = "Fruta"
$ws.Range("G403").Value = 100101
$ws.Range("H403").Value = "Berries"
$ws.Range("I403").Value = 100112025
$ws.Range("J403").Value = "Frutilla"
$ws.Range("K403").Value = "Sin especificar"
$ws.Range("L403").Value = "Segunda"
$ws.Range("M403").Value = 300
$ws.Range("N403").Value = 7500
$ws.Range("O403").Value = 8000
$ws.Range("P403").Value = 7750
$ws.Range("Q403").Value = "$/bandeja 7 kilos"
$ws.Range("R403").Value = "Provincia de Melipilla"
$ws.Range("S403").Value = 1107
$ws.Range("T403").Value = 7
